$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Containers")
$ws3.Columns.Item(1).ColumnWidth = 17
$ws3.Columns.Item(2).ColumnWidth = 10.28515625
$ws3.Columns.Item(3).ColumnWidth = 9.7109375
$ws3.Columns.Item(4).ColumnWidth = 18.7109375
$ws3.Columns.Item(5).ColumnWidth = 22
$ws3.Columns.Item(6).ColumnWidth = 15.7109375
$ws3.Columns.Item(7).ColumnWidth = 0
$ws3.Columns.Item(8).ColumnWidth = 1
$ws3.Columns.Item(9).ColumnWidth = 0.5
$ws3.Columns.Item(10).ColumnWidth = 2.5
Write-Host "done"
Write-Host "G1=" $ws3.Columns.Item(1).ColumnWidth
Write-Host "G2=" $ws3.Columns.Item(2).ColumnWidth
